$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO"
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M10").Value = 124.22
$ws1.Range("N10").Value = 2438.55
$ws1.Range("M22").Value = "2 de 20"

# Sheet "VENTA MENSUAL"
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F10").Value = 11565.71
$ws2.Range("F22").Value = 54669.1
